# ISRaD Kleber_2005.xlsx update: add "pro_usda_soil_order" field to the
# "profile" sheet (new column N) and its controlled-vocabulary list
# (new column E on the "controlled vocabulary" sheet), plus the dropdown
# data validation that ties them together.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. "controlled vocabulary" sheet: insert a new column E holding the
#    field name and the allowed USDA soil-order values.
# ---------------------------------------------------------------------
$cv = $wb.Worksheets.Item("controlled vocabulary")

$cv.Columns("E").Insert()

$cv.Range("E2").Value = "pro_usda_soil_order"

$cv.Range("E4").Value  = "Alfisols"
$cv.Range("E5").Value  = "Andisols"
$cv.Range("E6").Value  = "Aridisols"
$cv.Range("E7").Value  = "Entisols"
$cv.Range("E8").Value  = "Gelisols"
$cv.Range("E9").Value  = "Histosols"
$cv.Range("E10").Value = "Inceptisols"
$cv.Range("E11").Value = "Mollisols"
$cv.Range("E12").Value = "Oxisols"
$cv.Range("E13").Value = "Spodosols"
$cv.Range("E14").Value = "Ultisols"
$cv.Range("E15").Value = "Vertisols"

# ---------------------------------------------------------------------
# 2. "profile" sheet: insert a new column N ("pro_usda_soil_order")
#    and populate the USDA soil order for each profile.
# ---------------------------------------------------------------------
$pro = $wb.Worksheets.Item("profile")

$pro.Columns("N").Insert()

$pro.Range("N1").Value = "pro_usda_soil_order"

$pro.Range("N4").Value  = "Mollisols"
$pro.Range("N5").Value  = "Alfisols"
$pro.Range("N6").Value  = "Inceptisols"
$pro.Range("N7").Value  = "Inceptisols"
$pro.Range("N8").Value  = "Inceptisols"
$pro.Range("N9").Value  = "Spodosols"
$pro.Range("N10").Value = "Inceptisols"
$pro.Range("N11").Value = "Alfisols"
$pro.Range("N12").Value = "Andisols"
$pro.Range("N13").Value = "Alfisols"

# Dropdown list validation for the whole column, sourced from the
# controlled-vocabulary sheet's new list.
$valRange = $pro.Range("N4:N1048576")
$valRange.Validation.Add(3, 1, 1, "='controlled vocabulary'!`$E`$4:`$E`$15")

# ---------------------------------------------------------------------
# 3. "metadata" sheet: wrap the long citation text in M4 and let the
#    row grow to fit it.
# ---------------------------------------------------------------------
$meta = $wb.Worksheets.Item("metadata")
$meta.Range("M4").WrapText = $true
$meta.Rows(4).RowHeight = 409.6

# ---------------------------------------------------------------------
# 4. Restore/record the various sheet selections seen in the workbook
#    after the edit (cosmetic UI state).
# ---------------------------------------------------------------------
$cv.Range("D57").Select()
$pro.Range("O18").Select()

$layer = $wb.Worksheets.Item("layer")
$layer.Range("S57:X57").Select()

# "metadata" ends up as the active sheet/tab.
$meta.Activate()
$meta.Range("A4").Select()
